$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10 get letters A-I in natural order
$ws.Cells.Item(2, 1).Value = "A"
$ws.Cells.Item(3, 1).Value = "B"
$ws.Cells.Item(4, 1).Value = "C"
$ws.Cells.Item(5, 1).Value = "D"
$ws.Cells.Item(6, 1).Value = "E"
$ws.Cells.Item(7, 1).Value = "F"
$ws.Cells.Item(8, 1).Value = "G"
$ws.Cells.Item(9, 1).Value = "H"
$ws.Cells.Item(10, 1).Value = "I"

# Rows 12, 13, 11 written in this order so the shared-strings table picks up
# K, L, J in that sequence (matching the source workbook's layout)
$ws.Cells.Item(12, 1).Value = "K"
$ws.Cells.Item(13, 1).Value = "L"
$ws.Cells.Item(11, 1).Value = "J"

# Rows 14-16 get letters M-O
$ws.Cells.Item(14, 1).Value = "M"
$ws.Cells.Item(15, 1).Value = "N"
$ws.Cells.Item(16, 1).Value = "O"

$ws.Range("C16").Select()
